# Update the "Training Dashboard" worksheet:
#  - Column H (PERIOD TO EXPIRE) values each decrease by 1 (rows 3-22)
#  - Column I (LAST UPDATE) values change from "03-Nov-2025" to "04-Nov-2025" (rows 3-22)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Force column I to be treated as plain text so the date-like string is not
# auto-converted into a date serial number by Excel's type inference.
$ws.Range("I3:I22").NumberFormat = "@"

for ($row = 3; $row -le 22; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE
    $iCell.Value2 = "04-Nov-2025"
}
